$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C21").Value = 307
$ws.Range("D21").Value = 258
$ws.Range("E21").Value = 49
$ws.Range("F21").Value = 73.92550143266476
